$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

$ws.Range("C7").Value = "Terrain"
$ws.Range("D7").Value = "LF/Yann"
$ws.Range("E7").Value = "40/60"

$ws.Range("C8").Value = "Meteo"
$ws.Range("D8").Value = "LF/Yann"
$ws.Range("E8").Value = "60/40"

$ws.Range("C9").Value = "Jeu"
$ws.Range("D9").Value = "LF/Yann"

$ws.Range("D10").Value = "Yann"

$ws.Range("E9").Value = "20/80"

$ws.Range("C10").Value = "Test"

$ws.Range("C11").Value = "Rapport"
$ws.Range("D11").Value = "LF/Yann"
$ws.Range("E11").Value = "50/50"

$ws.Range("C12").Value = "Plante"
$ws.Range("D12").Value = "LF/Yann"
$ws.Range("E12").Value = "90/10"

$ws.Range("C13").Value = "Programme(potager,casepotager…)"
$ws.Range("D13").Value = "LF/Yann"
$ws.Range("E13").Value = "50/50"

$ws.Range("E13").Select()
